# Applies MAC address (C) / device-id (D) data and flips status (F) from
# "libre" to "en uso" for the rooms that now have an assigned access point.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row => MAC address (column C), device id (column D)
$newData = @{
    4 = @("3C:46:A1:25:40:30", "122379002563")
    5 = @("3C:46:A1:25:7E:40", "122379003345")
    6 = @("3C:46:A1:25:79:00", "122379003248")
    7 = @("3C:46:A1:25:7D:30", "122379003336")
    8 = @("3C:46:A1:25:3C:70", "122379002751")
    9 = @("3C:46:A1:25:3F:90", "122379002574")
    10 = @("3C:46:A1:25:3C:C0", "122379002492")
    11 = @($null, $null)
    12 = @($null, $null)
    13 = @($null, $null)
    14 = @($null, $null)
    15 = @($null, $null)
    16 = @($null, $null)
    17 = @($null, $null)
    18 = @($null, $null)
    19 = @($null, $null)
    20 = @($null, $null)
    21 = @($null, $null)
    22 = @($null, $null)
    23 = @($null, $null)
    24 = @($null, $null)
    25 = @($null, $null)
    26 = @($null, $null)
    27 = @("3C:46:A1:25:7B:60", "122379003491")
    28 = @("3C:46:A1:25:46:D0", "122379002766")
    29 = @("3C:46:A1:25:7B:30", "122379003496")
    30 = @("3C:46:A1:25:78:50", "122379003270")
    31 = @("3C:46:A1:25:7B:B0", "122379003479")
    32 = @("3C:46:A1:25:7C:B0", "122379003372")
    33 = @("3C:46:A1:25:47:00", "122379002655")
    34 = @("3C:46:A1:25:3C:50", "122379002753")
    35 = @("3C:46:A1:25:77:B0", "122379003291")
    36 = @("3C:46:A1:25:7C:A0", "122379003407")
    37 = @("3C:46:A1:25:48:00", "122379002661")
    38 = @("3C:46:A1:25:7B:40", "122379003495")
    39 = @("3C:46:A1:25:3C:A0", "122379002744")
    40 = @("3C:46:A1:25:45:C0", "122379002607")
    41 = @("3C:46:A1:25:79:10", "122379003247")
    42 = @("3C:46:A1:25:7B:D0", "122379003765")
    43 = @("3C:46:A1:25:3C:40", "122379002754")
    44 = @("3C:46:A1:25:44:90", "122379002624")
    45 = @("3C:46:A1:25:7D:C0", "122379003484")
    46 = @("3C:46:A1:25:77:A0", "122379003300")
    47 = @("3C:46:A1:25:7E:80", "122379003341")
}

foreach ($row in $newData.Keys) {
    $mac = $newData[$row][0]
    $devId = $newData[$row][1]

    if ($null -ne $mac) {
        $ws.Cells.Item($row, 3).Value = $mac
    }
    if ($null -ne $devId) {
        # Device id is all digits; force text storage so it keeps
        # leading context / matches the MAC-like id column, not a number.
        $ws.Cells.Item($row, 4).NumberFormat = "@"
        $ws.Cells.Item($row, 4).Value = $devId
    }

    # Every one of these rooms is now occupied.
    $ws.Cells.Item($row, 6).Value = "en uso"
}
